$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of row 2 and row 3 for columns D, M, N, O, P, R, S
$cols = @("D", "M", "N", "O", "P", "R", "S")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $tmp = $cell2.Value2
    $cell2.Value2 = $cell3.Value2
    $cell3.Value2 = $tmp
}
